# Autogenerated on Thu Mar 26 2015 18:06:15 GMT+0000 (Coordinated Universal Time)
#
# Rewrites the "Source:" footer block of the Summary sheet: the old two
# sources (book.moeasmea.gov.tw white paper + SMEA paragraph) are kept,
# but each is now followed by a blank spacer row, the old hyperlinked URL
# cell becomes a plain (non-hyperlinked) text cell, and a brand-new
# "Others:" source block (SMEA admin name + URL) plus a revised "SMEA"
# citation are appended further down the column.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-SourceCell($addr, $text) {
    # "source" named cell style = italic, non-bold, non-underline, default color
    $c = $ws.Range($addr)
    $c.Value = $text
    $c.Font.Bold = $false
    $c.Font.Italic = $true
    $c.Font.Underline = $false
    $c.Font.ColorIndex = -4105
}

function Set-TitleCell($addr, $text) {
    # "title" named cell style = bold, non-italic, non-underline, default color
    $c = $ws.Range($addr)
    $c.Value = $text
    $c.Font.Bold = $true
    $c.Font.Italic = $false
    $c.Font.Underline = $false
    $c.Font.ColorIndex = -4105
}

# Row 53 "Source:" is unchanged.

# Row 54 becomes a blank spacer (was the "White Paper..." text).
Set-SourceCell "A54" ""

# Row 55 now carries the "White Paper on SME..." text (used to be on row 54).
Set-SourceCell "A55" "White Paper on Small and Medium Enterprises in Taiwan, 2013"

# Row 56 stays a blank spacer.
Set-SourceCell "A56" ""

# Row 57 (new) carries the URL that used to live on A55 as a hyperlink;
# it is now plain text with the "source" style, no hyperlink.
Set-SourceCell "A57" "http://book.moeasmea.gov.tw/book/doc_detail.jsp?pub_SerialNo=2013A01165&click=2013A01165"

# Remove the old hyperlink (A55 -> rId1) entirely.
$ws.Hyperlinks.Delete()

# Row 58 (new) blank spacer.
Set-SourceCell "A58" ""

# Row 59: "Others:" label (was the bold "SMEA" title; now plain "source" style).
Set-SourceCell "A59" "Others:"

# Row 60 (new) blank spacer.
Set-SourceCell "A60" ""

# Row 61: new org-name text (replaces the old long "SMEA ... White Paper" text).
Set-SourceCell "A61" "Small and Medium Enterprise Administration - Ministry of Economic Affairs, Taiwan - China"

# Row 62 (new) blank spacer.
Set-SourceCell "A62" ""

# Row 63 (new): the "Others:" block's URL.
Set-SourceCell "A63" "http://www.moeasmea.gov.tw/mp.asp?mp=2"

# Row 66 (new): "SMEA" bold title, reappearing further down.
Set-TitleCell "A66" "SMEA"

# Row 67 (new): revised SMEA citation text.
Set-SourceCell "A67" "SMEA, revised and issued on September 2, 2009, available at http://www.moeasmea.gov.tw/ct.asp?xItem=70&CtNode=261&mp=2"

Write-Host "edit applied"
